$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 168
$ws.Range("F3").Value = 2431
$ws.Range("F6").Value = 81
$ws.Range("F7").Value = 293
$ws.Range("F8").Value = 361
$ws.Range("F9").Value = 3361
$ws.Range("F10").Value = 887
$ws.Range("F11").Value = 98
$ws.Range("F13").Value = 1536
$ws.Range("F14").Value = 10
$ws.Range("F15").Value = 888
$ws.Range("F16").Value = 1733
$ws.Range("F19").Value = 1514
$ws.Range("F21").Value = 83
$ws.Range("F23").Value = 4020
$ws.Range("F25").Value = 2652
$ws.Range("F26").Value = 1180
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 39
$ws.Range("F19").Value = 157
$ws.Range("F48").Value = 309
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2506
$ws.Range("F5").Value = 714
$ws.Range("F7").Value = 9557
$ws.Range("F12").Value = 2856
$ws.Range("F13").Value = 396
$ws.Range("F14").Value = 720
$ws.Range("F15").Value = 41
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2506
$ws.Range("F3").Value = 714
$ws.Range("F5").Value = 168
$ws.Range("F6").Value = 2856
$ws.Range("F7").Value = 397
$ws.Range("F9").Value = 720
$ws.Range("F10").Value = 39
$ws.Range("F14").Value = 81
$ws.Range("F15").Value = 293
$ws.Range("F16").Value = 361
$ws.Range("F18").Value = 98
$ws.Range("F21").Value = 10
$ws.Range("F24").Value = 888
$ws.Range("F27").Value = 1733
$ws.Range("F30").Value = 1514
$ws.Range("F36").Value = 83
$ws.Range("F41").Value = 4020
$ws.Range("F43").Value = 2652
$ws.Range("F47").Value = 309
$ws.Range("F48").Value = 1180
